# Weekly update: insert a new data row for "Feria Lagunitas de Puerto Montt -
# Ciboulette" just before the current row 223, pushing the rest of the table
# (previously rows 223-344) down by one row (to 224-345).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 223 - shifts rows 223:344 down to 224:345
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with this week's record
$ws.Cells.Item(223, 1).Value = 4
$ws.Cells.Item(223, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(223, 3).Value = "Los Lagos"
$ws.Cells.Item(223, 4).Value = 45029
$ws.Cells.Item(223, 5).Value = 10
$ws.Cells.Item(223, 6).Value = 100112039
$ws.Cells.Item(223, 7).Value = "Ciboulette"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 80
$ws.Cells.Item(223, 11).Value = 3500
$ws.Cells.Item(223, 12).Value = 3500
$ws.Cells.Item(223, 13).Value = 3500
$ws.Cells.Item(223, 14).Value = "`$/docena de atados"
$ws.Cells.Item(223, 15).Value = "Región Metropolitana"
$ws.Cells.Item(223, 16).Value = 1167
$ws.Cells.Item(223, 17).Value = 3
$ws.Cells.Item(223, 18).Value = "Hortaliza"
